# GestionRisques.xlsx - "ajustement gantt et gestion risque"
# Applies the changes described in the commit: widen column D, adjust a few
# row heights, update the responsable of risk row 4, populate two new risk
# rows (17 & 18) that used to be blank, and update the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D width -------------------------------------------------------
# Target stored width ~26.332; this engine stores width = round(chars*6)/6 + 5/6
# so chars = 25.5 lands the stored width on 26.333333... (closest achievable).
$ws.Columns.Item(4).ColumnWidth = 25.5

# --- Row height tweaks ------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 43.8
$ws.Rows.Item(13).RowHeight = 43.2
$ws.Rows.Item(14).RowHeight = 57.6

# --- Row 4: responsable changes from "M. Sylvain" to "E. Bourque" ----------
$ws.Range("I4").Value = "E. Bourque"

# --- Populate previously-blank rows 17 & 18 with new risk entries ----------
# First copy formatting from row 16 (the last populated data row) so the new
# rows inherit the same borders/alignment/number formats.
$ws.Range("B16:I16").Copy() | Out-Null
$ws.Range("B17:I17").PasteSpecial(-4122) | Out-Null
$ws.Range("B18:I18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(17).RowHeight = 72
$ws.Rows.Item(18).RowHeight = 72

# Row 17
$ws.Range("B17").Value = 'Non respect des délais de livraison des documents à produire pour les livrables'
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 'Nécessité par les autres membres de terminer le travail pour les autres ou devoir remettre un livrable incomplet'
$ws.Range("E17").Value = 'Perte de points liés aux sections manquantes ou moins raffinées du livrable'
$ws.Range("F17").Value = 0.15
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = 'Vérification périodique de l''avancement des travaux par le chef d''équipe ainsi que mise en place d''un "deadline" pour la remise des travaux un jour ou deux avant la remise du livrable'
$ws.Range("I17").Value = "M. Sylvain"

# Row 18 (H18's text is introduced before E18's in the source workbook, so
# we assign values in that same order to keep the shared-string table in the
# same append order as the authored workbook)
$ws.Range("B18").Value = 'Dépassement du budget du projet'
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 'Nécessité de retirer des pièces du robot pour arriver en dessous du 300$ imposé'
$ws.Range("F18").Value = 0.05
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = 'Assigner à un membre de l''équipe la charge de conserver un relevé du coût des pièces placées sur le robot ainsi que les factures'
$ws.Range("E18").Value = 'Performances réduites si le nouvel équipement de remplacement est moins performant ou plus complexe à intégrer'
$ws.Range("I18").Value = "M. Sylvain"

# --- View: scroll / selection ----------------------------------------------
# Best-effort; this headless runtime does not persist topLeftCell to the
# saved sheetView, but we still set the active selection to match.
$ws.Activate()
$ws.Range("E22").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
